$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextCell 'D2' '30.438.26'
Set-TextCell 'E2' '  +0.81%  '
Set-TextCell 'D3' '1.869.95'
Set-TextCell 'E3' '  +0.36%  '
Set-TextCell 'E4' '  -0.02%  '
Set-TextCell 'D5' '246.34'
Set-TextCell 'E5' '  +1.50%  '
Set-TextCell 'E6' '  +0.00%  '
Set-TextCell 'E7' '  +0.67%  '
Set-TextCell 'D8' '0.2914'
Set-TextCell 'E8' '  +2.20%  '
Set-TextCell 'D9' '0.06496'
Set-TextCell 'E9' '  +0.42%  '
Set-TextCell 'D10' '21.98'
Set-TextCell 'E10' '  +6.03%  '
Set-TextCell 'E11' '  +0.31%  '
Set-TextCell 'D12' '97.57'
Set-TextCell 'E12' '  +3.29%  '
Set-TextCell 'D13' '0.7393'
Set-TextCell 'E13' '  +8.45%  '
Set-TextCell 'D14' '1.873.25'
Set-TextCell 'E14' '  +0.42%  '
Set-TextCell 'D15' '5.121'
Set-TextCell 'E15' '  +0.96%  '
Set-TextCell 'D16' '274.63'
Set-TextCell 'E16' '  +2.29%  '
Set-TextCell 'D17' '30.410.46'
Set-TextCell 'E17' '  +0.72%  '
Set-TextCell 'D18' '13.37'
Set-TextCell 'E18' '  +0.35%  '
Set-TextCell 'D19' '0.000007550'
Set-TextCell 'E19' '  +0.30%  '
Set-TextCell 'D20' '1.0000'
Set-TextCell 'E20' '  +0.00%  '
Set-TextCell 'D21' '2.115.79'
Set-TextCell 'E21' '  +0.15%  '
Set-TextCell 'E22' '  -0.02%  '
Set-TextCell 'D23' '5.232'
Set-TextCell 'E23' '  +1.02%  '
Set-TextCell 'D24' '6.169'
Set-TextCell 'E24' '  +1.09%  '
Set-TextCell 'D25' '9.292'
Set-TextCell 'E25' '  -0.19%  '
Set-TextCell 'D26' '164.18'
Set-TextCell 'E26' '  -0.66%  '
Set-TextCell 'E27' '  +0.30%  '
Set-TextCell 'D28' '1.930'
Set-TextCell 'E28' '  +2.18%  '
Set-TextCell 'D29' '0.09986'
Set-TextCell 'E29' '  +1.81%  '
Set-TextCell 'D30' '1.366'
Set-TextCell 'E30' '  -0.47%  '
Set-TextCell 'D31' '1.506'
Set-TextCell 'E31' '  -0.36%  '
Set-TextCell 'D32' '4.304'
Set-TextCell 'E32' '  +1.76%  '
Set-TextCell 'D33' '4.147'
Set-TextCell 'E33' '  +4.41%  '
Set-TextCell 'D34' '0.04834'
Set-TextCell 'E34' '  +3.18%  '
Set-TextCell 'D35' '1.121'
Set-TextCell 'E35' '  +1.05%  '
Set-TextCell 'D36' '0.6969'
Set-TextCell 'E36' '  +1.64%  '
Set-TextCell 'B37' 'Frax'
Set-TextCell 'C37' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell 'D37' '0.9997'
Set-TextCell 'E37' '  +0.03%  '
Set-TextCell 'B38' 'HuobiToken'
Set-TextCell 'C38' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D38' '2.715'
Set-TextCell 'E38' '  +0.23%  '
Set-TextCell 'B39' 'VeChain'
Set-TextCell 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D39' '0.01858'
Set-TextCell 'E39' '  +0.83%  '
Set-TextCell 'B40' 'MXToken'
Set-TextCell 'C40' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D40' '2.734'
Set-TextCell 'E40' '  +0.04%  '
Set-TextCell 'B41' 'FraxShare'
Set-TextCell 'C41' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D41' '6.305'
Set-TextCell 'E41' '  -1.44%  '
Set-TextCell 'B42' 'Aave'
Set-TextCell 'C42' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D42' '72.79'
Set-TextCell 'E42' '  +3.73%  '
Set-TextCell 'B43' 'RenderToken'
Set-TextCell 'C43' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D43' '1.968'
Set-TextCell 'E43' '  +4.59%  '
Set-TextCell 'B44' 'TheSandbox'
Set-TextCell 'C44' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D44' '0.4191'
Set-TextCell 'E44' '  +3.23%  '
Set-TextCell 'B45' 'PaxDollar'
Set-TextCell 'C45' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D45' '0.9998'
Set-TextCell 'E45' '  -0.01%  '
Set-TextCell 'B46' 'TrustWalletToken'
Set-TextCell 'C46' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D46' '0.8339'
Set-TextCell 'E46' '  -0.41%  '
Set-TextCell 'B47' 'Quant'
Set-TextCell 'C47' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D47' '102.05'
Set-TextCell 'E47' '  +0.33%  '
Set-TextCell 'B48' 'EnergySwap'
Set-TextCell 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D48' '9.201'
Set-TextCell 'E48' '  +0.61%  '
Set-TextCell 'B49' 'Aptos'
Set-TextCell 'C49' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D49' '7.014'
Set-TextCell 'E49' '  +1.29%  '
Set-TextCell 'B50' 'Maker'
Set-TextCell 'C50' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 'D50' '929.35'
Set-TextCell 'E50' '  +1.28%  '
Set-TextCell 'B51' 'Elrond'
Set-TextCell 'C51' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 'D51' '35.30'
Set-TextCell 'E51' '  +2.61%  '
